$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6 (Hours, Minutes, Seconds, Miles, Pace)
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 16.1
$ws.Range("G2").Value = 19

$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 13.2
$ws.Range("G3").Value = 23

$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 17.5
$ws.Range("G4").Value = 21

$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 21.3
$ws.Range("G5").Value = 23

$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 15.6
$ws.Range("G6").Value = 28

# Add new rows 7-9
$ws.Range("A7").Value = 45297
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B7").Value = "trail"
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 42
$ws.Range("F7").Value = 15.6
$ws.Range("G7").Value = 35

$ws.Range("A8").Value = 45298
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B8").Value = "road"
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 23.6
$ws.Range("G8").Value = 15

$ws.Range("A9").Value = 45299
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B9").Value = "trail"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 25.9
$ws.Range("G9").Value = 14
